$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the two additional columns
$ws.Range("G1").Value = "yield_g"
$ws.Range("H1").Value = "harvest_date"

# New data for row 14 (plant_id 13 / Dig / E bed)
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 327.1

# H14 should look like the existing date column (E) -- copy its number
# format (style) over before writing the new value so it shares the same
# style index instead of minting a new one.
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("H14").Value = 45520

# Match column H's width to column E's (bestFit date column) as closely as
# possible.
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# Leave the selection where the author left off after entering the data.
$ws.Range("H15").Select() | Out-Null
